$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Asset_Name_Ex" column (C) for rows 4 and 5 -- append the
# "_UnwantedNameN" suffixes that were introduced in this revision.
$ws.Range("C4").Value = "a007Q00000CAA2f_ExcelNameCol_UnwantedName1"
$ws.Range("C5").Value = "a007Q00000CA9ne_ExcelNameCol_UnwantedName2"

# Update the "Asset_Description_Ex" column (E) for rows 3 and 4 -- append
# the "_UnwantedDescN" suffix to the shared description text.
$ws.Range("E3").Value = "Aluslevypari  M8, NL8SP_ExcelDescCol_UnwantedDesc1"
$ws.Range("E4").Value = "Aluslevypari  M8, NL8SP_ExcelDescCol_UnwantedDesc1"

# Column C's text got longer, so its best-fit width grows accordingly.
$ws.Columns.Item(3).ColumnWidth = 44.5

# Leave the selection on E5, matching the saved view state.
$ws.Range("E5").Select() | Out-Null
